# Auto-generated edit script: refresh market-price derived columns
# (currentAveragePrice / NQ / HQ, LevePrice NQ/HQ, LeveProfit NQ/HQ)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 100
$ws.Range("H100").Value = 22223366
$ws.Range("I100").Value = 22223366
$ws.Range("K100").Value = 22223366
$ws.Range("M100").Value = -22222825
# Row 123
$ws.Range("H123").Value = 42215
$ws.Range("J123").Value = 42215
$ws.Range("L123").Value = 42215
$ws.Range("N123").Value = -52015
# Row 132
$ws.Range("H132").Value = 99831.484
$ws.Range("I132").Value = 115757.03
$ws.Range("K132").Value = 347271.09
$ws.Range("M132").Value = -344741.09

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 557.26086
$ws.Range("I2").Value = 459.82352
$ws.Range("J2").Value = 833.3333
$ws.Range("K2").Value = 459.82352
$ws.Range("L2").Value = 833.3333
$ws.Range("M2").Value = -346.82352
$ws.Range("N2").Value = -1059.3333
# Row 3
$ws.Range("H3").Value = 7407.143
$ws.Range("J3").Value = 10990
$ws.Range("L3").Value = 10990
$ws.Range("N3").Value = -11220
# Row 32
$ws.Range("H32").Value = 7125.569
$ws.Range("I32").Value = 5326.161
$ws.Range("J32").Value = 9914.65
$ws.Range("K32").Value = 5326.161
$ws.Range("L32").Value = 9914.65
$ws.Range("M32").Value = -5039.161
$ws.Range("N32").Value = -10488.65
# Row 45
$ws.Range("H45").Value = 1146
$ws.Range("I45").Value = 1026.6666
$ws.Range("J45").Value = 1197.1428
$ws.Range("K45").Value = 1026.6666
$ws.Range("L45").Value = 1197.1428
$ws.Range("M45").Value = -649.6666
$ws.Range("N45").Value = -1951.1428
# Row 97
$ws.Range("H97").Value = 664
$ws.Range("I97").Value = 664
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 664
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -168
# Row 102
$ws.Range("H102").Value = 2500
$ws.Range("I102").Value = 2400
$ws.Range("K102").Value = 2400
$ws.Range("M102").Value = -778
# Row 116
$ws.Range("H116").Value = 557.26086
$ws.Range("I116").Value = 459.82352
$ws.Range("J116").Value = 833.3333
$ws.Range("K116").Value = 459.82352
$ws.Range("L116").Value = 833.3333
$ws.Range("M116").Value = 1834.17648
$ws.Range("N116").Value = -5421.3333
# Row 132
$ws.Range("H132").Value = 2372.509
$ws.Range("I132").Value = 1812.3182
$ws.Range("J132").Value = 4613.273
$ws.Range("K132").Value = 5436.9546
$ws.Range("L132").Value = 13839.819
$ws.Range("M132").Value = -2906.9546
$ws.Range("N132").Value = -18899.819

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 557.26086
$ws.Range("I3").Value = 459.82352
$ws.Range("J3").Value = 833.3333
$ws.Range("K3").Value = 459.82352
$ws.Range("L3").Value = 833.3333
$ws.Range("M3").Value = -345.82352
$ws.Range("N3").Value = -1061.3333
# Row 5
$ws.Range("H5").Value = 2183.4285
$ws.Range("J5").Value = 2826.25
$ws.Range("L5").Value = 2826.25
$ws.Range("N5").Value = -3052.25
# Row 99
$ws.Range("H99").Value = 2760.0667
$ws.Range("I99").Value = 1941.4286
$ws.Range("K99").Value = 1941.4286
$ws.Range("M99").Value = -443.4286

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 10872428
$ws.Range("I31").Value = 1530.4688
$ws.Range("K31").Value = 1530.4688
$ws.Range("M31").Value = -1235.4688
# Row 34
$ws.Range("H34").Value = 10872428
$ws.Range("I34").Value = 1530.4688
$ws.Range("K34").Value = 1530.4688
$ws.Range("M34").Value = -1328.4688
# Row 58
$ws.Range("H58").Value = 1791.0857
$ws.Range("I58").Value = 1558.678
$ws.Range("K58").Value = 1558.678
$ws.Range("M58").Value = -1355.678
# Row 82
$ws.Range("H82").Value = 39300
$ws.Range("J82").Value = 39300
$ws.Range("L82").Value = 39300
$ws.Range("N82").Value = -40022
# Row 85
$ws.Range("H85").Value = 39300
$ws.Range("J85").Value = 39300
$ws.Range("L85").Value = 39300
$ws.Range("N85").Value = -41796
# Row 132
$ws.Range("H132").Value = 2983.6875
$ws.Range("I132").Value = 1324.0667
$ws.Range("K132").Value = 3972.2001
$ws.Range("M132").Value = -1442.2001
# Row 136
$ws.Range("H136").Value = 1791.0857
$ws.Range("I136").Value = 1558.678
$ws.Range("K136").Value = 4676.034000000001
$ws.Range("M136").Value = -2126.034000000001

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 965814.0600000001
$ws.Range("J4").Value = 1627.1364
$ws.Range("L4").Value = 4881.4092
$ws.Range("N4").Value = -5105.4092
# Row 121
$ws.Range("H121").Value = 1689.375
$ws.Range("I121").Value = 247.8
$ws.Range("J121").Value = 1811.5424
$ws.Range("K121").Value = 743.4000000000001
$ws.Range("L121").Value = 5434.6272
$ws.Range("M121").Value = 566.5999999999999
$ws.Range("N121").Value = -8054.6272
# Row 140
$ws.Range("H140").Value = 3417.3333
$ws.Range("I140").Value = 3438.3333
$ws.Range("J140").Value = 3333.3333
$ws.Range("K140").Value = 10314.9999
$ws.Range("L140").Value = 9999.999899999999
$ws.Range("M140").Value = -5134.999899999999
$ws.Range("N140").Value = -20359.9999

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 55
$ws.Range("I2").Value = 41.42857
$ws.Range("J2").Value = 74
$ws.Range("K2").Value = 41.42857
$ws.Range("L2").Value = 74
$ws.Range("M2").Value = 71.57142999999999
$ws.Range("N2").Value = -300
# Row 4
$ws.Range("H4").Value = 28998
$ws.Range("J4").Value = 28998
$ws.Range("L4").Value = 28998
$ws.Range("N4").Value = -29222
# Row 70
$ws.Range("H70").Value = 5932.1147
$ws.Range("I70").Value = 5508.756
$ws.Range("J70").Value = 6800
$ws.Range("K70").Value = 5508.756
$ws.Range("L70").Value = 6800
$ws.Range("M70").Value = -5238.756
$ws.Range("N70").Value = -7340
# Row 73
$ws.Range("H73").Value = 5932.1147
$ws.Range("I73").Value = 5508.756
$ws.Range("J73").Value = 6800
$ws.Range("K73").Value = 5508.756
$ws.Range("L73").Value = 6800
$ws.Range("M73").Value = -4572.756
$ws.Range("N73").Value = -8672
# Row 135
$ws.Range("H135").Value = 48335.95
$ws.Range("J135").Value = 48335.95
$ws.Range("L135").Value = 48335.95
$ws.Range("N135").Value = -58475.95

$ws = $wb.Worksheets.Item("LTW")
# Row 81
$ws.Range("H81").Value = 58359.6
$ws.Range("J81").Value = 58359.6
$ws.Range("L81").Value = 58359.6
$ws.Range("N81").Value = -60355.6
# Row 84
$ws.Range("H84").Value = 58359.6
$ws.Range("J84").Value = 58359.6
$ws.Range("L84").Value = 175078.8
$ws.Range("N84").Value = -185062.8
# Row 94
$ws.Range("H94").Value = 31888.334
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 31888.334
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 31888.334
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -33240.334
# Row 132
$ws.Range("H132").Value = 2862.348
$ws.Range("I132").Value = 1616.7441
$ws.Range("J132").Value = 4922.385
$ws.Range("K132").Value = 4850.2323
$ws.Range("L132").Value = 14767.155
$ws.Range("M132").Value = -2320.2323
$ws.Range("N132").Value = -19827.155
# Row 138
$ws.Range("H138").Value = 59980
$ws.Range("J138").Value = 59980
$ws.Range("L138").Value = 59980
$ws.Range("N138").Value = -70260

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 11113149
$ws.Range("I132").Value = 1078.3529
$ws.Range("K132").Value = 3235.0587
$ws.Range("M132").Value = -705.0587000000005
# Row 136
$ws.Range("H136").Value = 1770.2727
$ws.Range("J136").Value = 2657.5
$ws.Range("L136").Value = 7972.5
$ws.Range("N136").Value = -13072.5
